$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh coin price / volume(1h) / name / link cells with the latest
# coinranking.com snapshot values. Cells in this sheet are plain text
# (numeric-looking prices and "x.xx%" volumes are stored as strings), so
# temporarily force a text number format before writing each value to stop
# Excel auto-converting it to a Number/Percentage, then clear the format
# again so the cell keeps the workbook's original (unstyled) look.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '310.85'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '0.83%'
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '39.29'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '1.95%'
$ws.Range('E3').ClearFormats()
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.166'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '1.25%'
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.08154'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '0.28%'
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.985'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '1.01%'
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '8.142'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '2.55%'
$ws.Range('E7').ClearFormats()
$ws.Range('B8').NumberFormat = "@"
$ws.Range('B8').Value = 'GateToken'
$ws.Range('B8').ClearFormats()
$ws.Range('C8').NumberFormat = "@"
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('C8').ClearFormats()
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '4.234'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '0.81%'
$ws.Range('E8').ClearFormats()
$ws.Range('B9').NumberFormat = "@"
$ws.Range('B9').Value = 'MXToken'
$ws.Range('B9').ClearFormats()
$ws.Range('C9').NumberFormat = "@"
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C9').ClearFormats()
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9264'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '-0.21%'
$ws.Range('E9').ClearFormats()
$ws.Range('B10').NumberFormat = "@"
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('B10').ClearFormats()
$ws.Range('C10').NumberFormat = "@"
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('C10').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1390'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-4.23%'
$ws.Range('E10').ClearFormats()
$ws.Range('B11').NumberFormat = "@"
$ws.Range('B11').Value = 'WazirX'
$ws.Range('B11').ClearFormats()
$ws.Range('C11').NumberFormat = "@"
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('C11').ClearFormats()
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.1930'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '-1.30%'
$ws.Range('E11').ClearFormats()
$ws.Range('B12').NumberFormat = "@"
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('B12').ClearFormats()
$ws.Range('C12').NumberFormat = "@"
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('C12').ClearFormats()
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.09028'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-1.04%'
$ws.Range('E12').ClearFormats()
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('B13').ClearFormats()
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('C13').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03509'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '0.22%'
$ws.Range('E13').ClearFormats()
$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('B14').ClearFormats()
$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('C14').ClearFormats()
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.09829'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '0.26%'
$ws.Range('E14').ClearFormats()
$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('B15').ClearFormats()
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('C15').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.001389'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '-0.76%'
$ws.Range('E15').ClearFormats()
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('B16').ClearFormats()
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('C16').ClearFormats()
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.005948'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '-2.27%'
$ws.Range('E16').ClearFormats()
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'LEO'
$ws.Range('B17').ClearFormats()
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C17').ClearFormats()
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.679'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '1.00%'
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.377'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '-2.24%'
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3456'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '0.28%'
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.1349'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '2.86%'
$ws.Range('E20').ClearFormats()
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.650'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '-3.17%'
$ws.Range('E21').ClearFormats()
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-1.27%'
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04365'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '-1.43%'
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.001229'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '0.87%'
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.004871'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '0.49%'
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0001300'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-0.08%'
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0003996'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '-10.15%'
$ws.Range('E27').ClearFormats()
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02142'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '1.23%'
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.05203'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '1.56%'
$ws.Range('E40').ClearFormats()
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007428'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-0.48%'
$ws.Range('E41').ClearFormats()
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.009836'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '-3.02%'
$ws.Range('E42').ClearFormats()
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1367'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '0.20%'
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.002130'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-0.54%'
$ws.Range('E44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.009887'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '5.05%'
$ws.Range('E45').ClearFormats()
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00006384'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '2.64%'
$ws.Range('E46').ClearFormats()
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000749'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-0.17%'
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0009991'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '-37.59%'
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.002685'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '-12.47%'
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00002098'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '-0.17%'
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0001998'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '-0.17%'
$ws.Range('E51').ClearFormats()
